# StaticData 외래키 테스트 (#73)
# Add a new "ForeignTest" worksheet (placed after "GroupedItemTest", i.e. as
# the last tab) with a small Id/TargetTestId/Value/StudentId/비고 table that
# demonstrates a foreign-key-style join where the related student rows may
# or may not be present.

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "ForeignTest"
$newSheet.Move($null, $wb.Worksheets("GroupedItemTest"))

$ws = $wb.Worksheets("ForeignTest")

# Title cell pointing at where the real table starts (same convention used
# by the other sheets in this workbook, e.g. "A7", "B2").
$ws.Range("A1").Value = "C9"

# Header row
$ws.Range("C9").Value = "Id"
$ws.Range("D9").Value = "TargetTestId"
$ws.Range("E9").Value = "Value"
$ws.Range("F9").Value = "StudentId"
$ws.Range("G9").Value = "비고"

# Data rows
$ws.Range("C10").Value = 1001
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = "AAA"
$ws.Range("F10").Value = 20220001
$ws.Range("G10").Value = "학생이 있을수도 있고,"

$ws.Range("C11").Value = 1002
$ws.Range("D11").Value = 102
$ws.Range("E11").Value = "BBB"
$ws.Range("G11").Value = "학생이 없을수도 있습니다."

$ws.Range("C12").Value = 1003
$ws.Range("D12").Value = 104
$ws.Range("E12").Value = "CCC"
$ws.Range("F12").Value = 20220002

# Make the new sheet the active tab/selection, like the captured workbook.
$ws.Activate()
$ws.Range("H22").Select() | Out-Null
